# CSV 데이터를 Scriptable Object로 자동 변환
# Insert a new "Name" (internal key) column between Id and DisplayName in
# the research table, shifting DisplayName..RequireResearchId one column
# to the right (B..F -> C..G).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the comment that currently documents ResearchAmount (C1) and
# the table before we start moving things around.
$lo = $ws.ListObjects.Item(1)
$cmt = $ws.Range("C1").Comment
$cmtText = $cmt.Text()
$lo.Unlist()
$cmt.Delete()

# Shift the existing DisplayName..RequireResearchId columns one slot to
# the right (rightmost column first so values are not clobbered), using
# plain value copies rather than a structural column insert so column
# widths/metadata are left untouched.
for ($r = 1; $r -le 9; $r++) {
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 6).Value()
}
for ($r = 1; $r -le 9; $r++) {
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 5).Value()
}
for ($r = 1; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 4).Value()
}
for ($r = 1; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Value()
}
for ($r = 1; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 2).Value()
}

# Fill the now-freed column B with the new "Name" header + internal keys,
# one per research row (same order as DisplayName).
$ws.Range("B1").Value = "Name"
$ws.Range("B2").Value = "Copper_Tech"
$ws.Range("B3").Value = "Copper_UpProduce"
$ws.Range("B4").Value = "Copper_Discount"
$ws.Range("B5").Value = "Copper_Drill"
$ws.Range("B6").Value = "Iron_Tech"
$ws.Range("B7").Value = "Iron_UpProduce"
$ws.Range("B8").Value = "Iron_Discount"
$ws.Range("B9").Value = "Iron_Drill"

# Rebuild the table over the new A1:G9 range so it picks up the headers
# (Id, Name, DisplayName, ResearchAmount, InputItemPerTickId,
# InputItemPerTickCount, RequireResearchId) in the correct order.
$newlo = $ws.ListObjects.Add(1, $ws.Range("A1:G9"), $null, 1)
$newlo.Name = "표1"

# The comment that documented ResearchAmount now belongs on D1, its new
# location after the column shift.
$ws.Range("D1").AddComment($cmtText)

# Restore the previously active selection, shifted to its new location.
$ws.Range("C8").Select()
